$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.597.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.967.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.844"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.255.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.966.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.524.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  +7.22%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +20.61%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.92%  "
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  -12.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0971"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.370.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.150.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.96%  "
